# Setting up ChattyISS.pptx - reword agenda bullet + highlight key filenames
$p = $ppt.ActivePresentation

# --- Slide 2 ("Agenda"): reword bullet text ---------------------------------
$slide2 = $p.Slides.Item(2)
$agendaBox = $slide2.Shapes.Item(2)
$agendaTr = $agendaBox.TextFrame.TextRange
$agendaTr.Paragraphs(2, 1).Text = "Migration procedure"

# --- Slide 5 ("Setting up"): highlight the 5 python filenames in green -----
$slide5 = $p.Slides.Item(5)
$filesBox = $slide5.Shapes.Item(2)
$filesTr = $filesBox.TextFrame.TextRange

# Paragraphs 2-6 are: chatbot.py, chatbot_gpt4.py, imagebot.py, ingestion.py,
# understandingimages_llama.py. Build the combined character range covering
# exactly those paragraphs (inclusive) and apply a bright-green highlight.
$firstPara = $filesTr.Paragraphs(2, 1)
$lastPara = $filesTr.Paragraphs(6, 1)
$startChar = $firstPara.Start
$endChar = $lastPara.Start + $lastPara.Length
$length = $endChar - $startChar

$pyFilesRange = $filesTr.Characters($startChar, $length)
$pyFilesRange.Font.Highlight = 65280  # RGB(0, 255, 0) -> 0x00FF00, bright green
